$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 549, shifting existing rows 549:614 down to 550:615.
$ws.Rows("549:549").Insert()

# Populate the newly-inserted row 549 with its data (same categorical fields as
# the row that follows it, but its own measurement values).
$ws.Range("A549").Value = 9
$ws.Range("B549").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C549").Value = "Metropolitana"
$ws.Range("D549").Value = 44858
$ws.Range("E549").Value = 13
$ws.Range("F549").Value = 100112031
$ws.Range("G549").Value = "Poroto verde"
$ws.Range("H549").Value = "Magnum"
$ws.Range("I549").Value = "Primera"
$ws.Range("J549").Value = 45
$ws.Range("K549").Value = 38000
$ws.Range("L549").Value = 38000
$ws.Range("M549").Value = 38000
$ws.Range("N549").Value = "$/malla 25 kilos"
$ws.Range("O549").Value = "Perú"
$ws.Range("P549").Value = 1520
$ws.Range("Q549").Value = 25
$ws.Range("R549").Value = "Hortaliza"
